# Applies the "output generated at 456a3b4" update to 苏州-漫展信息.xlsx
# Touches two sheets that carry duplicated data: 展览 (sheet1) and 全部类型 (sheet4).
# For each sheet:
#   - a handful of "想去人数" (F column) counters are bumped up
#   - three consecutive rows (an ONLY event rotates out, two newer events rotate in)
#     get their name/venue/time/stats/link/cover updated

$wb = $excel.ActiveWorkbook

function Set-SimpleCounts {
    param($ws, $map)
    foreach ($row in $map.Keys) {
        $ws.Cells.Item([int]$row, 6).Value2 = $map[$row]
    }
}

function Set-TextCell {
    # Writes a plain text value into a cell without letting Excel's
    # date auto-detection reinterpret date-shaped strings (e.g. 2024-05-18)
    # as date serial numbers.
    param($ws, $row, $col, $text)
    $looksLikeIsoDate = $text -match '^\d{4}-\d{2}-\d{2}$'
    if ($looksLikeIsoDate) {
        $ws.Cells.Item($row, $col).NumberFormat = "@"
    }
    $ws.Cells.Item($row, $col).Value2 = $text
}

function Update-EventRow {
    param($ws, $row, $fields)
    if ($fields.ContainsKey("B")) { Set-TextCell $ws $row 2 $fields["B"] }
    if ($fields.ContainsKey("C")) { Set-TextCell $ws $row 3 $fields["C"] }
    if ($fields.ContainsKey("D")) { Set-TextCell $ws $row 4 $fields["D"] }
    if ($fields.ContainsKey("E")) { Set-TextCell $ws $row 5 $fields["E"] }
    if ($fields.ContainsKey("F")) { $ws.Cells.Item($row, 6).Value2 = $fields["F"] }
    if ($fields.ContainsKey("G")) { $ws.Cells.Item($row, 7).Value2 = $fields["G"] }
    if ($fields.ContainsKey("H")) { Set-TextCell $ws $row 8 $fields["H"] }
    if ($fields.ContainsKey("I")) { Set-TextCell $ws $row 9 $fields["I"] }
}

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")

Set-SimpleCounts $ws1 @{
    2  = 621
    4  = 1290
    5  = 1162
    7  = 16462
    10 = 10
    22 = 135
    26 = 6681
    32 = 5742
    36 = 4802
}

Update-EventRow $ws1 28 @{
    B = "2024-05-18"
    C = "苏州·OrangeOrange国潮&随机宅舞派对【免费活动】"
    D = "狮山路298号 金鹰国际购物中心(狮山路店)"
    E = "2024.05.18 13:00-05.18 17:00"
    F = 18
    G = 29
    H = "https://show.bilibili.com/platform/detail.html?id=83949"
    I = "//i1.hdslb.com/bfs/openplatform/202404/DOH6BK8i1712638105049.png"
}

Update-EventRow $ws1 29 @{
    C = "苏州·YoungComic动漫嘉年华"
    D = "清禾路886号 尹山湖大剧院"
    E = "2024.05.18 10:00-05.18 17:00"
    F = 1116
    G = 60
    H = "https://show.bilibili.com/platform/detail.html?id=83142"
    I = "//i2.hdslb.com/bfs/openplatform/202403/4wWLK6Jg1710840463319.jpeg"
}

Update-EventRow $ws1 30 @{
    C = "苏州·明日方舟ONLY#2024~佑桑柔"
    D = "城际路21号 苏州汇融广场假日酒店"
    F = 10
    G = 75
    H = "https://show.bilibili.com/platform/detail.html?id=84046"
    I = "//i0.hdslb.com/bfs/openplatform/202404/t4T75Yi31712890052782.jpeg"
}

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")

Set-SimpleCounts $ws4 @{
    2  = 621
    4  = 1290
    5  = 1162
    7  = 16462
    10 = 10
    22 = 135
    27 = 6681
    35 = 5742
    39 = 4802
}

Update-EventRow $ws4 29 @{
    B = "2024-05-18"
    C = "苏州·OrangeOrange国潮&随机宅舞派对【免费活动】"
    D = "狮山路298号 金鹰国际购物中心(狮山路店)"
    E = "2024.05.18 13:00-05.18 17:00"
    F = 18
    G = 29
    H = "https://show.bilibili.com/platform/detail.html?id=83949"
    I = "//i1.hdslb.com/bfs/openplatform/202404/DOH6BK8i1712638105049.png"
}

Update-EventRow $ws4 30 @{
    C = "苏州·YoungComic动漫嘉年华"
    D = "清禾路886号 尹山湖大剧院"
    E = "2024.05.18 10:00-05.18 17:00"
    F = 1116
    G = 60
    H = "https://show.bilibili.com/platform/detail.html?id=83142"
    I = "//i2.hdslb.com/bfs/openplatform/202403/4wWLK6Jg1710840463319.jpeg"
}

Update-EventRow $ws4 31 @{
    C = "苏州·明日方舟ONLY#2024~佑桑柔"
    D = "城际路21号 苏州汇融广场假日酒店"
    F = 10
    G = 75
    H = "https://show.bilibili.com/platform/detail.html?id=84046"
    I = "//i0.hdslb.com/bfs/openplatform/202404/t4T75Yi31712890052782.jpeg"
}

Write-Host "Edit complete"
